# Weekly update: a new price record (week of 2023-03-28) is inserted for
# "Feria Lagunitas de Puerto Montt - Pepino dulce" right after the existing
# row 76, pushing the previously-last rows (77-84) down by one (to 78-85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 77:84 down to 78:85 by inserting a new blank row at 77.
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new weekly record.
$ws.Range("A77").Value = 4
$ws.Range("B77").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C77").Value = "Los Lagos"
$ws.Range("D77").Value = 45013
$ws.Range("E77").Value = 10
$ws.Range("F77").Value = 100112043
$ws.Range("G77").Value = "Pepino dulce"
$ws.Range("H77").Value = "Cultivar IV Región"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 80
$ws.Range("K77").Value = 17000
$ws.Range("L77").Value = 18000
$ws.Range("M77").Value = 17500
$ws.Range("N77").Value = "`$/bandeja 18 kilos"
$ws.Range("O77").Value = "Provincia de Limarí"
$ws.Range("P77").Value = 972
$ws.Range("Q77").Value = 18
$ws.Range("R77").Value = "Hortaliza"

# Make sure the new row's date cell keeps the same date/time style used by
# the rest of column D.
$ws.Range("D77").NumberFormat = $ws.Range("D78").NumberFormat
